$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 8564.429
$ws.Range("I74").Value = 7289.2144
$ws.Range("J74").Value = 11114.857
$ws.Range("K74").Value = 7289.2144
$ws.Range("L74").Value = 11114.857
$ws.Range("M74").Value = -6353.2144
$ws.Range("N74").Value = -12986.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 8564.429
$ws.Range("I77").Value = 7289.2144
$ws.Range("J77").Value = 11114.857
$ws.Range("K77").Value = 36446.072
$ws.Range("L77").Value = 55574.285
$ws.Range("M77").Value = -31766.072
$ws.Range("N77").Value = -64934.285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1634.8572
$ws.Range("I101").Value = 1298.6666
$ws.Range("J101").Value = 2240
$ws.Range("K101").Value = 3895.9998
$ws.Range("L101").Value = 6720
$ws.Range("M101").Value = -2273.9998
$ws.Range("N101").Value = -9964

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 5577.0527
$ws.Range("I131").Value = 4458.846
$ws.Range("J131").Value = 7999.8335
$ws.Range("K131").Value = 13376.538
$ws.Range("L131").Value = 23999.5005
$ws.Range("M131").Value = -8336.537999999999
$ws.Range("N131").Value = -34079.50049999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 76929200
$ws.Range("I137").Value = 142860860
$ws.Range("J137").Value = 8934
$ws.Range("K137").Value = 428582580
$ws.Range("L137").Value = 26802
$ws.Range("M137").Value = -428580030
$ws.Range("N137").Value = -31902

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6075.2666
$ws.Range("I138").Value = 3934.2727
$ws.Range("J138").Value = 6767.9414
$ws.Range("K138").Value = 11802.8181
$ws.Range("L138").Value = 20303.8242
$ws.Range("M138").Value = -6662.8181
$ws.Range("N138").Value = -30583.8242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1584.2561
$ws.Range("I32").Value = 1035.6533
$ws.Range("J32").Value = 7462.143
$ws.Range("K32").Value = 1035.6533
$ws.Range("L32").Value = 7462.143
$ws.Range("M32").Value = -748.6532999999999
$ws.Range("N32").Value = -8036.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 111114744
$ws.Range("I45").Value = 200001810
$ws.Range("J45").Value = 5912.25
$ws.Range("K45").Value = 200001810
$ws.Range("L45").Value = 5912.25
$ws.Range("M45").Value = -200001433
$ws.Range("N45").Value = -6666.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1891.0303
$ws.Range("I110").Value = 1427.0358
$ws.Range("J110").Value = 4489.4
$ws.Range("K110").Value = 1427.0358
$ws.Range("L110").Value = 4489.4
$ws.Range("M110").Value = 617.9641999999999
$ws.Range("N110").Value = -8579.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6761.625
$ws.Range("I132").Value = 3156
$ws.Range("J132").Value = 7963.5
$ws.Range("K132").Value = 9468
$ws.Range("L132").Value = 23890.5
$ws.Range("M132").Value = -6938
$ws.Range("N132").Value = -28950.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6280.5713
$ws.Range("I134").Value = 1962.5
$ws.Range("K134").Value = 5887.5
$ws.Range("M134").Value = -3352.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 178.35
$ws.Range("I7").Value = 59.090908
$ws.Range("J7").Value = 324.1111
$ws.Range("K7").Value = 59.090908
$ws.Range("L7").Value = 324.1111
$ws.Range("M7").Value = 53.909092
$ws.Range("N7").Value = -550.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1107.6957
$ws.Range("I22").Value = 540.7646999999999
$ws.Range("J22").Value = 2714
$ws.Range("K22").Value = 540.7646999999999
$ws.Range("L22").Value = 2714
$ws.Range("M22").Value = -190.7646999999999
$ws.Range("N22").Value = -3414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33355.566
$ws.Range("I31").Value = 3824.8262
$ws.Range("J31").Value = 81870.36
$ws.Range("K31").Value = 3824.8262
$ws.Range("L31").Value = 81870.36
$ws.Range("M31").Value = -3529.8262
$ws.Range("N31").Value = -82460.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 33355.566
$ws.Range("I34").Value = 3824.8262
$ws.Range("J34").Value = 81870.36
$ws.Range("K34").Value = 3824.8262
$ws.Range("L34").Value = 81870.36
$ws.Range("M34").Value = -3622.8262
$ws.Range("N34").Value = -82274.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4085.1177
$ws.Range("I132").Value = 3515.2693
$ws.Range("J132").Value = 5937.125
$ws.Range("K132").Value = 10545.8079
$ws.Range("L132").Value = 17811.375
$ws.Range("M132").Value = -8015.8079
$ws.Range("N132").Value = -22871.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2607.6858
$ws.Range("I134").Value = 2089.3125
$ws.Range("J134").Value = 8137
$ws.Range("K134").Value = 6267.9375
$ws.Range("L134").Value = 24411
$ws.Range("M134").Value = -3732.9375
$ws.Range("N134").Value = -29481

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1505.5555
$ws.Range("I14").Value = 1505.5555
$ws.Range("K14").Value = 4516.666499999999
$ws.Range("M14").Value = -4343.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 298
$ws.Range("J92").Value = 298
$ws.Range("L92").Value = 894
$ws.Range("N92").Value = -3390

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10236380
$ws.Range("J131").Value = 27780348
$ws.Range("L131").Value = 83341044
$ws.Range("N131").Value = -83351124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10115
$ws.Range("I80").Value = 9000
$ws.Range("J80").Value = 10561
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 10561
$ws.Range("M80").Value = -8002
$ws.Range("N80").Value = -12557

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10115
$ws.Range("I83").Value = 9000
$ws.Range("J83").Value = 10561
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 52805
$ws.Range("M83").Value = -40008
$ws.Range("N83").Value = -62789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4137.7856
$ws.Range("J68").Value = 8158.8
$ws.Range("L68").Value = 8158.8
$ws.Range("N68").Value = -9656.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4137.7856
$ws.Range("J71").Value = 8158.8
$ws.Range("L71").Value = 40794
$ws.Range("N71").Value = -48282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4715.1665
$ws.Range("I100").Value = 2898.4119
$ws.Range("J100").Value = 9127.286
$ws.Range("K100").Value = 2898.4119
$ws.Range("L100").Value = 9127.286
$ws.Range("M100").Value = -2357.4119
$ws.Range("N100").Value = -10209.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3337.7693
$ws.Range("I132").Value = 1693.2667
$ws.Range("J132").Value = 8819.444
$ws.Range("K132").Value = 5079.800099999999
$ws.Range("L132").Value = 26458.332
$ws.Range("M132").Value = -2549.800099999999
$ws.Range("N132").Value = -31518.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6521.2
$ws.Range("I136").Value = 3203.7778
$ws.Range("K136").Value = 9611.3334
$ws.Range("M136").Value = -7061.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3766.3333
$ws.Range("I122").Value = 3030.9473
$ws.Range("J122").Value = 10752.5
$ws.Range("K122").Value = 9092.841899999999
$ws.Range("L122").Value = 32257.5
$ws.Range("M122").Value = -6642.841899999999
$ws.Range("N122").Value = -37157.5
